# Weekly update: add the newest "Fruta / Piña" price record for
# Vega Monumental Concepción at the top of the data block (row 264),
# pushing the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 264; this shifts rows
# 264:309 down to 265:310 (and the sheet's used range grows to T310).
$ws.Rows.Item(264).EntireRow.Insert()

# Populate the new row with this week's record.
$ws.Range("A264").Value = 11
$ws.Range("B264").Value = "Vega Monumental Concepción"
$ws.Range("C264").Value = "Bíobío"
$ws.Range("D264").Value = 45211
$ws.Range("E264").Value = 8
$ws.Range("F264").Value = "Fruta"
$ws.Range("G264").Value = 100108
$ws.Range("H264").Value = "Tropicales y subtropicales"
$ws.Range("I264").Value = 100108005
$ws.Range("J264").Value = "Piña"
$ws.Range("K264").Value = "Caramelo"
$ws.Range("L264").Value = "Segunda"
$ws.Range("M264").Value = 200
$ws.Range("N264").Value = 20000
$ws.Range("O264").Value = 21000
$ws.Range("P264").Value = 20500
$ws.Range("Q264").Value = "`$/caja 14 unidades"
$ws.Range("R264").Value = "Ecuador"
$ws.Range("S264").Value = 1464
$ws.Range("T264").Value = 14
